# "Generate Report for Handoff"
#
# The localization status report moved from "In Translation" to
# "Ready for handoff" and the handoff/xliff-generation timestamps were
# refreshed. The same status string and the same "Latest HO Xliff
# Generate Date" / "Latest Handoff Datetime" timestamp are shared
# across the Overview sheet and the per-locale (zh-cn / de-de) sheets,
# so update every cell that surfaces that data.
#
# NOTE: these date/time cells are stored as literal text (not real
# Excel date serials) even though a date/time number format is applied
# to them, so plain string assignment via .Value keeps their cell type
# intact.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview!E2 (zh-cn status), Overview!F2 (de-de status)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
# Overview!G2 "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-08-26 04:37:51"

# zh-cn!C2 "Status"
$wsZhCn.Range("C2").Value = "Ready for handoff"
# zh-cn!H2 "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-08-26 04:37:47"

# de-de!C2 "Status"
$wsDeDe.Range("C2").Value = "Ready for handoff"
# de-de!H2 "Latest Handoff Datetime"
$wsDeDe.Range("H2").Value = "2016-08-26 04:37:51"

# "Ready for handoff" is longer than "In Translation", so Excel widened
# the Status columns that display it: Overview columns E/F and column C
# on the zh-cn / de-de sheets.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
